$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 305). The whole column was bumped forward by one day,
# from serial 45203 (2023-10-04) to serial 45204 (2023-10-05).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 305 }

$ws.Range("C2:C$lastRow").Value = 45204
